# Auto-generated edit script: updates market-price-derived columns (H-N)
# on multiple Leve-profit worksheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 346
$ws.Range("J41").Value = 439.5
$ws.Range("L41").Value = 439.5
$ws.Range("N41").Value = -1319.5

$ws.Range("H80").Value = 10594844
$ws.Range("I80").Value = 346.14285
$ws.Range("J80").Value = 15229936
$ws.Range("K80").Value = 1038.42855
$ws.Range("L80").Value = 45689808
$ws.Range("M80").Value = -40.42855000000009
$ws.Range("N80").Value = -45691804

$ws.Range("H83").Value = 10594844
$ws.Range("I83").Value = 346.14285
$ws.Range("J83").Value = 15229936
$ws.Range("K83").Value = 3115.28565
$ws.Range("L83").Value = 137069424
$ws.Range("M83").Value = 1876.71435
$ws.Range("N83").Value = -137079408

$ws.Range("H98").Value = 562.7273
$ws.Range("I98").Value = 611.25
$ws.Range("J98").Value = 433.33334
$ws.Range("K98").Value = 611.25
$ws.Range("L98").Value = 433.33334
$ws.Range("M98").Value = 886.75
$ws.Range("N98").Value = -3429.33334

$ws.Range("H106").Value = 6946850
$ws.Range("I106").Value = 12347333
$ws.Range("J106").Value = 3372.0952
$ws.Range("K106").Value = 12347333
$ws.Range("L106").Value = 3372.0952
$ws.Range("M106").Value = -12346702
$ws.Range("N106").Value = -4634.0952

$ws.Range("H112").Value = 3004002
$ws.Range("J112").Value = 3473269.5
$ws.Range("L112").Value = 10419808.5
$ws.Range("N112").Value = -10422024.5

$ws.Range("H122").Value = 562.7273
$ws.Range("I122").Value = 611.25
$ws.Range("J122").Value = 433.33334
$ws.Range("K122").Value = 1833.75
$ws.Range("L122").Value = 1300.00002
$ws.Range("M122").Value = 616.25
$ws.Range("N122").Value = -6200.000019999999

$ws.Range("H129").Value = 173453.72
$ws.Range("J129").Value = 182889.38
$ws.Range("L129").Value = 548668.14
$ws.Range("N129").Value = -558668.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6227.926
$ws.Range("I32").Value = 5385.88
$ws.Range("K32").Value = 5385.88
$ws.Range("M32").Value = -5098.88

$ws.Range("H88").Value = 126411.625
$ws.Range("I88").Value = 1248.4
$ws.Range("K88").Value = 1248.4
$ws.Range("M88").Value = -842.4000000000001

$ws.Range("H91").Value = 126411.625
$ws.Range("I91").Value = 1248.4
$ws.Range("K91").Value = 1248.4
$ws.Range("M91").Value = 155.5999999999999

$ws.Range("H102").Value = 2078.5
$ws.Range("I102").Value = 1063.3334
$ws.Range("J102").Value = 2513.5715
$ws.Range("K102").Value = 1063.3334
$ws.Range("L102").Value = 2513.5715
$ws.Range("M102").Value = 558.6666
$ws.Range("N102").Value = -5757.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1576.3513
$ws.Range("I86").Value = 1421.875
$ws.Range("J86").Value = 1861.5385
$ws.Range("K86").Value = 1421.875
$ws.Range("L86").Value = 1861.5385
$ws.Range("M86").Value = -298.875
$ws.Range("N86").Value = -4107.538500000001

$ws.Range("H89").Value = 1576.3513
$ws.Range("I89").Value = 1421.875
$ws.Range("J89").Value = 1861.5385
$ws.Range("K89").Value = 7109.375
$ws.Range("L89").Value = 9307.692500000001
$ws.Range("M89").Value = -1493.375
$ws.Range("N89").Value = -20539.6925

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 10000
$ws.Range("I103").Value = 10000
$ws.Range("K103").Value = 10000
$ws.Range("M103").Value = -8828

$ws.Range("H132").Value = 3118.52
$ws.Range("I132").Value = 2368.8125
$ws.Range("J132").Value = 4451.3335
$ws.Range("K132").Value = 7106.4375
$ws.Range("L132").Value = 13354.0005
$ws.Range("M132").Value = -4576.4375
$ws.Range("N132").Value = -18414.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1301.4286
$ws.Range("I68").Value = 1200
$ws.Range("K68").Value = 3600
$ws.Range("M68").Value = -2789

$ws.Range("H71").Value = 1301.4286
$ws.Range("I71").Value = 1200
$ws.Range("K71").Value = 10800
$ws.Range("M71").Value = -6744

$ws.Range("H131").Value = 742.53
$ws.Range("J131").Value = 744.11224
$ws.Range("L131").Value = 2232.33672
$ws.Range("N131").Value = -12312.33672

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").Value = ""

$ws.Range("H58").Value = 16672833
$ws.Range("J58").Value = 20006000
$ws.Range("L58").Value = 20006000
$ws.Range("N58").Value = -20006554

$ws.Range("H95").Value = 21562.666
$ws.Range("J95").Value = 21562.666
$ws.Range("L95").Value = 21562.666
$ws.Range("N95").Value = -27054.666

$ws.Range("H132").Value = 32482.588
$ws.Range("I132").Value = 3775
$ws.Range("J132").Value = 48141.273
$ws.Range("K132").Value = 11325
$ws.Range("L132").Value = 144423.819
$ws.Range("M132").Value = -8795
$ws.Range("N132").Value = -149483.819

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""

$ws.Range("H135").Value = 39770
$ws.Range("J135").Value = 39770
$ws.Range("L135").Value = 39770
$ws.Range("N135").Value = -49910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3401.8333
$ws.Range("I40").Value = 2833.3684
$ws.Range("J40").Value = 5562
$ws.Range("K40").Value = 2833.3684
$ws.Range("L40").Value = 5562
$ws.Range("M40").Value = -2697.3684
$ws.Range("N40").Value = -5834

$ws.Range("H61").Value = 4721.2144
$ws.Range("I61").Value = 1624.625
$ws.Range("J61").Value = 8850
$ws.Range("K61").Value = 1624.625
$ws.Range("L61").Value = 8850
$ws.Range("M61").Value = -1422.625
$ws.Range("N61").Value = -9254

$ws.Range("H68").Value = 2609
$ws.Range("J68").Value = 2827.8572
$ws.Range("L68").Value = 2827.8572
$ws.Range("N68").Value = -4325.8572

$ws.Range("H69").Value = 30000
$ws.Range("J69").Value = 30000
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622

$ws.Range("H71").Value = 2609
$ws.Range("J71").Value = 2827.8572
$ws.Range("L71").Value = 14139.286
$ws.Range("N71").Value = -21627.286

$ws.Range("H72").Value = 30000
$ws.Range("J72").Value = 30000
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112

$ws.Range("H113").Value = 4721.2144
$ws.Range("I113").Value = 1624.625
$ws.Range("J113").Value = 8850
$ws.Range("K113").Value = 1624.625
$ws.Range("L113").Value = 8850
$ws.Range("M113").Value = 545.375
$ws.Range("N113").Value = -13190

$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H136").Value = 1484
$ws.Range("I136").Value = 1571.3846
$ws.Range("K136").Value = 4714.1538
$ws.Range("M136").Value = -2164.1538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1620.381
$ws.Range("I132").Value = 1060.5883
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 3181.7649
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -651.7648999999997
$ws.Range("N132").Value = -17058.5
